$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Differentiate the three "Columbia Hatchery, 2016" source rows (13-15) with
# A/B/C suffixes, and replace the placeholder "Na" text in column K with the
# actual measured sodium values.
$ws.Range("A13").Value = "Columbia Hatchery, A, 2016"
$ws.Range("K13").Value = 1450

$ws.Range("A14").Value = "Columbia Hatchery, B, 2016"
$ws.Range("K14").Value = 1400

$ws.Range("A15").Value = "Columbia Hatchery, C, 2016"
$ws.Range("K15").Value = 1010

# Those three rows had an explicit custom row height; restore them to the
# sheet's standard (auto-fit) height.
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Rows.Item(14).EntireRow.AutoFit()
$ws.Rows.Item(15).EntireRow.AutoFit()

# Update the active cell / selection on the sheet.
[void]$ws.Range("H19").Select()
